# Add contoh (example) rows to the template worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New example data rows
$ws.Range("A2").Value = "IN000"
$ws.Range("B2").Value = "Contoh Mata Kuliah 1"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 72

$ws.Range("A3").Value = "IN001"
$ws.Range("B3").Value = "Contoh Mata Kuliah 2"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 72

$ws.Range("A4").Value = "IN002"
$ws.Range("B4").Value = "Contoh Mata Kuliah 3"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 72

# Column widths (no longer best-fit, custom widths). The COM layer here
# snaps ColumnWidth to a pixel grid (like real Excel), so these inputs are
# chosen to land on the stored width closest to the template's target.
$ws.Columns.Item(1).ColumnWidth = 27.42
$ws.Columns.Item(2).ColumnWidth = 31.1667
$ws.Columns.Item(3).ColumnWidth = 11.8333
$ws.Columns.Item(4).ColumnWidth = 22.6667

# Selection moves to B6
$ws.Range("B6").Select()

$wb.Save()
